# profile_df_2250.xlsx: add MinRollCrn / MaxRollCrn columns (D, E) to the
# crlc "Crn" vector profile, per "change imple of wr_grn_cr_scalar and fix
# mistakes in crlc Crn(including vector)".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers in row 1.
$ws.Cells.Item(1, 4).Value = "MinRollCrn"
$ws.Cells.Item(1, 5).Value = "MaxRollCrn"

# New min/max roll-crown values for data rows 2-8 (rows 2-5 use -0.9,
# rows 6-8 use -0.5; every row's max is 0.3).
$minRollCrn = @(-0.9, -0.9, -0.9, -0.9, -0.5, -0.5, -0.5)
$maxRollCrn = @(0.3, 0.3, 0.3, 0.3, 0.3, 0.3, 0.3)

for ($i = 0; $i -lt 7; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $minRollCrn[$i]
    $ws.Cells.Item($row, 5).Value = $maxRollCrn[$i]
}

# Leave the sheet with the same active selection captured in the edit.
$ws.Range("J8").Select()
